$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 463.732605
$ws.Range("H2").Value = 1391.197815
$ws.Range("I2").Value = 0.3632113435366598
$ws.Range("J2").Value = 0.3632113435366598
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 98.946724
$ws.Range("N2").Value = 296.840172
$ws.Range("O2").Value = 0.2098009692989996
$ws.Range("P2").Value = 0.2098009692989996
$ws.Range("Q2").Value = 45884.82207673602
$ws.Range("R2").Value = 412963.3986906242
$ws.Range("S2").Value = 0.07620209193438317
$ws.Range("T2").Value = 0.07620209193438317

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 463.732605
$ws.Range("H3").Value = 1391.197815
$ws.Range("I3").Value = 0.3632113435366598
$ws.Range("J3").Value = 0.3632113435366598
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 163.0062356666667
$ws.Range("N3").Value = 489.018707
$ws.Range("O3").Value = 0.345629090707923
$ws.Range("P3").Value = 0.3456290907079231
$ws.Range("Q3").Value = 75591.30629694724
$ws.Range("R3").Value = 680321.7566725252
$ws.Range("S3").Value = 0.1255364064013788
$ws.Range("T3").Value = 0.1255364064013788

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 463.732605
$ws.Range("H4").Value = 1391.197815
$ws.Range("I4").Value = 0.3632113435366598
$ws.Range("J4").Value = 0.3632113435366598
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 65.39610666666668
$ws.Range("N4").Value = 196.18832
$ws.Range("O4").Value = 0.1386621609326595
$ws.Range("P4").Value = 0.1386621609326595
$ws.Range("Q4").Value = 30326.30690139121
$ws.Range("R4").Value = 272936.7621125209
$ws.Range("S4").Value = 0.05036366977004779
$ws.Range("T4").Value = 0.0503636697700478

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 463.732605
$ws.Range("H5").Value = 1391.197815
$ws.Range("I5").Value = 0.3632113435366598
$ws.Range("J5").Value = 0.3632113435366598
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 144.2727966666667
$ws.Range("N5").Value = 432.81839
$ws.Range("O5").Value = 0.3059077790604178
$ws.Range("P5").Value = 0.3059077790604179
$ws.Range("Q5").Value = 66903.99982886865
$ws.Range("R5").Value = 602135.9984598178
$ws.Range("S5").Value = 0.1111091754308501
$ws.Range("T5").Value = 0.1111091754308501

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 12.24662533333333
$ws.Range("H6").Value = 36.739876
$ws.Range("I6").Value = 0.009591978638444229
$ws.Range("J6").Value = 0.009591978638444227
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 98.946724
$ws.Range("N6").Value = 296.840172
$ws.Range("O6").Value = 0.2098009692989996
$ws.Range("P6").Value = 0.2098009692989996
$ws.Range("Q6").Value = 1211.763456788741
$ws.Range("R6").Value = 10905.87111109867
$ws.Range("S6").Value = 0.002012406415840897
$ws.Range("T6").Value = 0.002012406415840897

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.24662533333333
$ws.Range("H7").Value = 36.739876
$ws.Range("I7").Value = 0.009591978638444229
$ws.Range("J7").Value = 0.009591978638444227
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 163.0062356666667
$ws.Range("N7").Value = 489.018707
$ws.Range("O7").Value = 0.345629090707923
$ws.Range("P7").Value = 0.3456290907079231
$ws.Range("Q7").Value = 1996.276295206704
$ws.Range("R7").Value = 17966.48665686033
$ws.Range("S7").Value = 0.0033152668548953
$ws.Range("T7").Value = 0.0033152668548953

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 12.24662533333333
$ws.Range("H8").Value = 36.739876
$ws.Range("I8").Value = 0.009591978638444229
$ws.Range("J8").Value = 0.009591978638444227
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 65.39610666666668
$ws.Range("N8").Value = 196.18832
$ws.Range("O8").Value = 0.1386621609326595
$ws.Range("P8").Value = 0.1386621609326595
$ws.Range("Q8").Value = 800.8816166053691
$ws.Range("R8").Value = 7207.934549448321
$ws.Range("S8").Value = 0.001330044485626585
$ws.Range("T8").Value = 0.001330044485626585

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 12.24662533333333
$ws.Range("H9").Value = 36.739876
$ws.Range("I9").Value = 0.009591978638444229
$ws.Range("J9").Value = 0.009591978638444227
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 144.2727966666667
$ws.Range("N9").Value = 432.81839
$ws.Range("O9").Value = 0.3059077790604178
$ws.Range("P9").Value = 0.3059077790604179
$ws.Range("Q9").Value = 1766.854886568849
$ws.Range("R9").Value = 15901.69397911964
$ws.Range("S9").Value = 0.002934260882081445
$ws.Range("T9").Value = 0.002934260882081445

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 689.7685036666667
$ws.Range("H10").Value = 2069.305511
$ws.Range("I10").Value = 0.5402504422695089
$ws.Range("J10").Value = 0.5402504422695089
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 98.946724
$ws.Range("N10").Value = 296.840172
$ws.Range("O10").Value = 0.2098009692989996
$ws.Range("P10").Value = 0.2098009692989996
$ws.Range("Q10").Value = 68250.33375619866
$ws.Range("R10").Value = 614253.0038057879
$ws.Range("S10").Value = 0.1133450664523562
$ws.Range("T10").Value = 0.1133450664523562

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 689.7685036666667
$ws.Range("H11").Value = 2069.305511
$ws.Range("I11").Value = 0.5402504422695089
$ws.Range("J11").Value = 0.5402504422695089
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 163.0062356666667
$ws.Range("N11").Value = 489.018707
$ws.Range("O11").Value = 0.345629090707923
$ws.Range("P11").Value = 0.3456290907079231
$ws.Range("Q11").Value = 112436.5672641327
$ws.Range("R11").Value = 1011929.105377194
$ws.Range("S11").Value = 0.1867262691161636
$ws.Range("T11").Value = 0.1867262691161637

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 689.7685036666667
$ws.Range("H12").Value = 2069.305511
$ws.Range("I12").Value = 0.5402504422695089
$ws.Range("J12").Value = 0.5402504422695089
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 65.39610666666668
$ws.Range("N12").Value = 196.18832
$ws.Range("O12").Value = 0.1386621609326595
$ws.Range("P12").Value = 0.1386621609326595
$ws.Range("Q12").Value = 45108.1746410924
$ws.Range("R12").Value = 405973.5717698316
$ws.Range("S12").Value = 0.07491229376991509
$ws.Range("T12").Value = 0.07491229376991511

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 689.7685036666667
$ws.Range("H13").Value = 2069.305511
$ws.Range("I13").Value = 0.5402504422695089
$ws.Range("J13").Value = 0.5402504422695089
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 144.2727966666667
$ws.Range("N13").Value = 432.81839
$ws.Range("O13").Value = 0.3059077790604178
$ws.Range("P13").Value = 0.3059077790604179
$ws.Range("Q13").Value = 99514.83107657192
$ws.Range("R13").Value = 895633.4796891473
$ws.Range("S13").Value = 0.165266812931074
$ws.Range("T13").Value = 0.165266812931074

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 111.00921
$ws.Range("H14").Value = 333.02763
$ws.Range("I14").Value = 0.08694623555538696
$ws.Range("J14").Value = 0.08694623555538696
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 98.946724
$ws.Range("N14").Value = 296.840172
$ws.Range("O14").Value = 0.2098009692989996
$ws.Range("P14").Value = 0.2098009692989996
$ws.Range("Q14").Value = 10983.99766332804
$ws.Range("R14").Value = 98855.97896995235
$ws.Range("S14").Value = 0.01824140449641933
$ws.Range("T14").Value = 0.01824140449641933

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 111.00921
$ws.Range("H15").Value = 333.02763
$ws.Range("I15").Value = 0.08694623555538696
$ws.Range("J15").Value = 0.08694623555538696
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 163.0062356666667
$ws.Range("N15").Value = 489.018707
$ws.Range("O15").Value = 0.345629090707923
$ws.Range("P15").Value = 0.3456290907079231
$ws.Range("Q15").Value = 18095.19344643049
$ws.Range("R15").Value = 162856.7410178744
$ws.Range("S15").Value = 0.03005114833548528
$ws.Range("T15").Value = 0.03005114833548529

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 111.00921
$ws.Range("H16").Value = 333.02763
$ws.Range("I16").Value = 0.08694623555538696
$ws.Range("J16").Value = 0.08694623555538696
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 65.39610666666668
$ws.Range("N16").Value = 196.18832
$ws.Range("O16").Value = 0.1386621609326595
$ws.Range("P16").Value = 0.1386621609326595
$ws.Range("Q16").Value = 7259.570138142401
$ws.Range("R16").Value = 65336.13124328161
$ws.Range("S16").Value = 0.01205615290706998
$ws.Range("T16").Value = 0.01205615290706999

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 111.00921
$ws.Range("H17").Value = 333.02763
$ws.Range("I17").Value = 0.08694623555538696
$ws.Range("J17").Value = 0.08694623555538696
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 144.2727966666667
$ws.Range("N17").Value = 432.81839
$ws.Range("O17").Value = 0.3059077790604178
$ws.Range("P17").Value = 0.3059077790604179
$ws.Range("Q17").Value = 16015.6091824573
$ws.Range("R17").Value = 144140.4826421157
$ws.Range("S17").Value = 0.02659752981641236
$ws.Range("T17").Value = 0.02659752981641236
